# FC_Gallery_Accessories.xlsx edit script
# Adds Denmark, Sweden, Norway and Italy market sheets (copied from the UK
# sheet's layout/format) and updates the active-sheet selection.

$wb = $excel.ActiveWorkbook

$uk = $wb.Worksheets.Item("UK")

# ---------------------------------------------------------------------
# Helper: build a "short" market sheet (Denmark / Sweden / Norway style)
# by copying UK, trimming the accessory list down to a single FCXB-S row
# and filling in the market name / part number.
# ---------------------------------------------------------------------
function New-ShortMarketSheet($afterSheet, $name, $marketText, $partNumber) {
    $uk.Copy([System.Reflection.Missing]::Value, $afterSheet)
    $ws = $wb.ActiveSheet
    $ws.Name = $name

    # Remove the ANC1...FC-ANC-E rows (original rows 11-18), leaving the
    # PSU rows above and the Wg/Accessories rows to close up beneath them.
    $ws.Range("A11:A18").EntireRow.Delete()

    # Re-open a single row for the FCXB-S accessory line and copy the
    # formatting from the row that will sit right below it (Wg).
    $ws.Range("A11").EntireRow.Insert()
    $ws.Range("A12").Copy()
    $ws.Range("A11").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("B2").Value2 = $marketText
    $ws.Range("B4").Value2 = $partNumber
    $ws.Range("A11").Value2 = "FCXB-S"

    return $ws
}

# ---------------------------------------------------------------------
# Denmark
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark = New-ShortMarketSheet $lastSheet "Denmark" "Denmark Market" "NGC-3446/T2004"

# ---------------------------------------------------------------------
# Italy (created now so that, after the later re-ordering, its sheetId
# ends up between Denmark's and Sweden's, matching the source history)
# ---------------------------------------------------------------------
$uk.Copy([System.Reflection.Missing]::Value, $denmark)
$italy = $wb.ActiveSheet
$italy.Name = "Italy"
$italy.Range("A17:A18").EntireRow.Delete()
$italy.Range("B4").ClearContents()

# ---------------------------------------------------------------------
# Sweden
# ---------------------------------------------------------------------
$sweden = New-ShortMarketSheet $italy "Sweden" "Sweden market" "NGC-3465/T2029"

# ---------------------------------------------------------------------
# Norway
# ---------------------------------------------------------------------
$norway = New-ShortMarketSheet $sweden "Norway" "Norway market" "NGC-3464/T1918"

# Move Italy so the final tab order is UK, Belgium, Denmark, Sweden, Norway, Italy.
# Re-fetch sheet references by name afterwards - stale object handles captured
# before a Move can end up pointing at the wrong tab.
$italyRef = $wb.Worksheets.Item("Italy")
$norwayRef = $wb.Worksheets.Item("Norway")
$italyRef.Move([System.Reflection.Missing]::Value, $norwayRef)

# Restore each sheet's selection / cursor position to match the target state.
$wb.Worksheets.Item("Denmark").Range("A1:XFD1048576").Select()
$wb.Worksheets.Item("Sweden").Range("A1:XFD1048576").Select()
$wb.Worksheets.Item("Italy").Range("A11").Select()

# Norway is the sheet left selected/active in the saved workbook.
$norwayFinal = $wb.Worksheets.Item("Norway")
$norwayFinal.Select()
$norwayFinal.Range("A7").Select()
